$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns B..G on the two ticket-data rows (2 and 3) hold numeric ticket
# counts (processing/completed/canceled/deferred/closed/new_or_reopened).
# Append the ":formatN()" filter to the placeholder text and switch the
# cell number format from Text ("@") to Number ("0") so the exported
# values render as numbers.
$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in @(2, 3)) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $text = $cell.Value2
        if ($text -like "*:formatN()*") {
            continue
        }
        $newText = $text -replace "\}$", ":formatN()}"
        $cell.Value = $newText
        $cell.NumberFormat = "0"
    }
}
